$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# row, D (Price), E (Volume 1h), G (Hora) -- $null means "no change for this cell"
$updates = @(
    @{ Row = 2; D = "312.82"; E = "5.32%"; G = "21" }
    @{ Row = 3; D = "44.31"; E = "7.03%"; G = "21" }
    @{ Row = 4; D = "5.104"; E = "1.22%"; G = "21" }
    @{ Row = 5; D = "0.07989"; E = "5.90%"; G = "21" }
    @{ Row = 6; D = "4.466"; E = "1.96%"; G = "21" }
    @{ Row = 7; D = "1.637"; E = "3.18%"; G = "21" }
    @{ Row = 8; D = "1.075"; E = "15.97%"; G = "21" }
    @{ Row = 9; D = "0.1286"; E = "6.14%"; G = "21" }
    @{ Row = 10; D = "0.1890"; E = "3.77%"; G = "21" }
    @{ Row = 11; D = "0.09210"; E = "3.69%"; G = "21" }
    @{ Row = 12; D = "0.04221"; E = "4.89%"; G = "21" }
    @{ Row = 13; D = $null; E = "-1.81%"; G = "21" }
    @{ Row = 14; D = "0.001308"; E = "1.70%"; G = "21" }
    @{ Row = 15; D = "0.005704"; E = "-4.32%"; G = "21" }
    @{ Row = 16; D = $null; E = $null; G = "21" }
    @{ Row = 17; D = "3.375"; E = "1.05%"; G = "21" }
    @{ Row = 18; D = "2.401"; E = "-0.90%"; G = "21" }
    @{ Row = 19; D = "0.3359"; E = "1.20%"; G = "21" }
    @{ Row = 20; D = "7.996"; E = "0.17%"; G = "21" }
    @{ Row = 21; D = "0.1375"; E = "-3.26%"; G = "21" }
    @{ Row = 22; D = "0.3123"; E = "3.97%"; G = "21" }
    @{ Row = 23; D = "0.04183"; E = "3.12%"; G = "21" }
    @{ Row = 24; D = "0.001270"; E = "0.30%"; G = "21" }
    @{ Row = 25; D = "0.004582"; E = "13.83%"; G = "21" }
    @{ Row = 26; D = "0.0001335"; E = "8.43%"; G = "21" }
    @{ Row = 27; D = $null; E = $null; G = "21" }
    @{ Row = 28; D = $null; E = $null; G = "21" }
    @{ Row = 29; D = $null; E = $null; G = "21" }
    @{ Row = 30; D = $null; E = $null; G = "21" }
    @{ Row = 31; D = $null; E = $null; G = "21" }
    @{ Row = 32; D = $null; E = $null; G = "21" }
    @{ Row = 33; D = $null; E = $null; G = "21" }
    @{ Row = 34; D = $null; E = $null; G = "21" }
    @{ Row = 35; D = $null; E = $null; G = "21" }
    @{ Row = 36; D = $null; E = $null; G = "21" }
    @{ Row = 37; D = $null; E = $null; G = "21" }
    @{ Row = 38; D = "0.02657"; E = "9.95%"; G = "21" }
    @{ Row = 39; D = "0.05416"; E = "3.91%"; G = "21" }
    @{ Row = 40; D = "0.005604"; E = "-14.52%"; G = "21" }
    @{ Row = 41; D = "0.007736"; E = "-0.80%"; G = "21" }
    @{ Row = 42; D = "0.1410"; E = "6.03%"; G = "21" }
    @{ Row = 43; D = "0.007292"; E = "-3.63%"; G = "21" }
    @{ Row = 44; D = "0.008385"; E = "6.88%"; G = "21" }
    @{ Row = 45; D = "0.3110"; E = "-3.42%"; G = "21" }
    @{ Row = 46; D = "0.00006714"; E = "-1.07%"; G = "21" }
    @{ Row = 47; D = "0.00000000742"; E = "-1.22%"; G = "21" }
    @{ Row = 48; D = "0.06199"; E = "34.60%"; G = "21" }
    @{ Row = 49; D = "0.003957"; E = "-5.93%"; G = "21" }
    @{ Row = 50; D = "0.00002077"; E = "-1.22%"; G = "21" }
    @{ Row = 51; D = "0.0001978"; E = "-1.22%"; G = "21" }
)

foreach ($u in $updates) {
    if ($null -ne $u.D) {
        $cell = $ws.Range("D" + $u.Row)
        $cell.NumberFormat = "@"
        $cell.Value = $u.D
    }
    if ($null -ne $u.E) {
        $cell = $ws.Range("E" + $u.Row)
        $cell.NumberFormat = "@"
        $cell.Value = $u.E
    }
    if ($null -ne $u.G) {
        $cell = $ws.Range("G" + $u.Row)
        $cell.NumberFormat = "@"
        $cell.Value = $u.G
    }
}
